# Fix PER bug: the Team (column B) / PER value (column C) pairs were
# mis-aligned. Re-point each row at the correct team code and replace the
# PER numeric value with the corrected figure.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value  = "POR"
$ws.Range("C2").Value  = 15.18461538461538

$ws.Range("B3").Value  = "NJN"
$ws.Range("C3").Value  = 11.16428571428571

$ws.Range("B4").Value  = "CLE"
$ws.Range("C4").Value  = 14.33076923076923

$ws.Range("B5").Value  = "DAL"
$ws.Range("C5").Value  = 15.42222222222222

$ws.Range("B6").Value  = "ATL"
$ws.Range("C6").Value  = 13.18571428571429

$ws.Range("B7").Value  = "OKC"
$ws.Range("C7").Value  = 13.14375

$ws.Range("B8").Value  = "CHA"
$ws.Range("C8").Value  = 12.47272727272727

$ws.Range("B9").Value  = "WAS"
$ws.Range("C9").Value  = 14.97272727272727

$ws.Range("B10").Value = "MIL"
$ws.Range("C10").Value = 13.08333333333333

$ws.Range("B11").Value = "LAC"
$ws.Range("C11").Value = 10.09166666666667

$ws.Range("B12").Value = "SAS"
$ws.Range("C12").Value = 15.55384615384615

$ws.Range("B13").Value = "DET"
$ws.Range("C13").Value = 13.03571428571429

$ws.Range("B14").Value = "ORL"
$ws.Range("C14").Value = 15.40833333333334

$ws.Range("B15").Value = "UTA"
$ws.Range("C15").Value = 14.46153846153846

$ws.Range("B16").Value = "MEM"
$ws.Range("C16").Value = 14.15714285714286

$ws.Range("B17").Value = "HOU"
$ws.Range("C17").Value = 11.12727272727273

$ws.Range("B18").Value = "DEN"
$ws.Range("C18").Value = 13.60769230769231

$ws.Range("B19").Value = "LAL"
$ws.Range("C19").Value = 13.36923076923077

$ws.Range("B20").Value = "GSW"
$ws.Range("C20").Value = 15.1

$ws.Range("B21").Value = "IND"
$ws.Range("C21").Value = 13.3

$ws.Range("B22").Value = "CHI"
$ws.Range("C22").Value = 10.66363636363636

$ws.Range("B23").Value = "PHI"
$ws.Range("C23").Value = 14.26363636363636

$ws.Range("B24").Value = "BOS"
$ws.Range("C24").Value = 13.81666666666666

$ws.Range("B25").Value = "TOR"
$ws.Range("C25").Value = 14.46428571428572

$ws.Range("B26").Value = "MIA"
$ws.Range("C26").Value = 12.24615384615385

$ws.Range("B27").Value = "SAC"
$ws.Range("C27").Value = 12.25384615384615

$ws.Range("B28").Value = "PHO"
$ws.Range("C28").Value = 13.07142857142857

$ws.Range("B29").Value = "NOH"
$ws.Range("C29").Value = 14.16363636363636

$ws.Range("B30").Value = "NYK"
$ws.Range("C30").Value = 12.96666666666667

$ws.Range("B31").Value = "MIN"
$ws.Range("C31").Value = 12.09230769230769
